# Update automatico via Actualizar 06-15-2020 15-08-31
#
# This script mirrors a routine "refresh the consultation dates / publish
# dates" pass over the DATACOVID "trabajo" sheet: the "Fecha consulta"
# (H) and "Fecha publicacion" (I) columns get re-stamped with current
# values and a cleaner dd-mm-yy date format, a couple of individual rows
# get their own updated dates, and the sheet selection is left on the
# range that was just touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Apply the new "dd-mm-yy" custom date format to the date columns ----
# (both H "Fecha consulta" and I "Fecha publicacion", rows 2-23)
$ws.Range("H2:I23").NumberFormat = "dd\-mm\-yy;@"

# --- Column H ("Fecha consulta") -----------------------------------------
# Rows 2-20 all get re-stamped to the same new consultation date.
$ws.Range("H2:H20").Value = 43987

# Row 21 gets its own refreshed consultation date.
$ws.Range("H21").Value = 43990

# Rows 22-23 keep their text placeholder value ("13-06-2020") - only the
# number format above changes for them.

# --- Column I ("Fecha publicacion") ---------------------------------------
# Most rows keep their existing publish date/text; only rows 21-23 change.
$ws.Range("I21").Value = 43989
$ws.Range("I22").Value = 43994
$ws.Range("I23").Value = 43993

# --- Sheet view: scroll over a column and leave the H2:H9 block selected --
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H2:H9").Select()
